$d = $word.ActiveDocument

$replacements = @(
    @("2025-02-22 Saturday", "2025-02-23 Sunday"),
    @("11×96=", "15×19="),
    @("84×66=", "50×79="),
    @("76×36=", "57×75="),
    @("51×41=", "91×44="),
    @("17×92=", "56×24="),
    @("75×59=", "86×73="),
    @("15×73=", "39×58="),
    @("58×45=", "71×77="),
    @("49×53=", "55×34="),
    @("24×97=", "69×65="),
    @("44×81=", "16×88="),
    @("46×48=", "22×21="),
    @("14×19=", "32×98="),
    @("89×38=", "48×68="),
    @("15×68=", "66×15="),
    @("39×90=", "98×16="),
    @("40×60=", "94×61="),
    @("43×32=", "66×72="),
    @("25×67=", "68×27="),
    @("99×97=", "45×15="),
    @("20×79=", "69×99="),
    @("49×93=", "51×15="),
    @("53×17=", "35×95="),
    @("72×18=", "54×13="),
    @("34×67=", "58×49=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
